# "Plot delay vs BPT" - reduce the sheet to a single column (A) of delay
# values against BPT, dropping the other series that used to live in
# columns B:D and updating the remaining values in column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop columns B:D entirely - only the BPT/delay column (A) remains.
$ws.Range("B1:D6").Delete()

# Update the values that remain in column A with the new delay series.
$ws.Range("A2").Value = 7226.425080128205
$ws.Range("A3").Value = 7.99
$ws.Range("A4").Value = 250
$ws.Range("A5").Value = 1
$ws.Range("A6").Value = 0.0001996662892080171

# Match the page margins Excel wrote for the cleaned-up sheet.
$ws.PageSetup.LeftMargin = 0.7 * 72
$ws.PageSetup.RightMargin = 0.7 * 72
$ws.PageSetup.TopMargin = 0.75 * 72
$ws.PageSetup.BottomMargin = 0.75 * 72
$ws.PageSetup.HeaderMargin = 0.3 * 72
$ws.PageSetup.FooterMargin = 0.3 * 72
